$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.920.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'1.703.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").Value = "'315.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "'0.4062"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("D8").Value = "'0.4067"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'53.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("D11").Value = "'1.469"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("D12").Value = "'0.08822"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'25.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("D14").Value = "'7.527"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "'8.061"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "'0.00001353"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'1.702.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "'96.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("D19").Value = "'0.07180"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "'21.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").Value = "'7.246"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "'14.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "'24.915.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'2.327"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").Value = "'6.792"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +29.87%  "
$ws.Range("D27").Value = "'2.892"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.87%  "
$ws.Range("D28").Value = "'23.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "'165.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'145.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("D31").Value = "'8.280"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.48%  "
$ws.Range("D32").Value = "'2.271"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.23%  "
$ws.Range("D33").Value = "'1.896.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "'0.08787"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'0.03206"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.39%  "
$ws.Range("D36").Value = "'7.327"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.51%  "
$ws.Range("D37").Value = "'1.017"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.2850"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("D39").Value = "'0.8495"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.99%  "
$ws.Range("D40").Value = "'10.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").Value = "'0.09414"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("D43").Value = "'17.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.19%  "
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "'2.718"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("D46").Value = "'0.7450"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("D47").Value = "'4.242"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'1.392"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("D49").Value = "'1.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'142.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").Value = "'0.08362"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.88%  "